# Ajusta na calculadora de IPCA
# Adds three new rows (31-33) of investment data to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=31; A="MAXIMA"; C=46059; D="CDB"; E="PRE";  F=0.10920000000000001; G="Easynvest"; H=43539; I=2000 },
    @{ Row=32; A="MAXIMA"; C=46479; D="CDB"; E="PRE";  F=0.11310000000000001; G="Easynvest"; H=43559; I=1000 },
    @{ Row=33; A="MAXIMA"; C=44697; D="CDB"; E="IPCA"; F=0.053; G="Easynvest"; H=43601; I=3000 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $cA = $ws.Cells.Item($row, 1)
    $cA.Value2 = $r.A
    $cA.Font.Name = "Calibri"
    $cA.Font.Size = 10

    $cB = $ws.Cells.Item($row, 2)
    $cB.Formula = "=YEAR(IF(C$row=""Liq. Diária"", NOW(), C$row))"
    $cB.Font.Name = "Calibri"
    $cB.Font.Size = 10

    $cC = $ws.Cells.Item($row, 3)
    $cC.Value2 = $r.C
    $cC.NumberFormat = "dd/mm/yyyy"
    $cC.Font.Name = "Calibri"
    $cC.Font.Size = 10

    $cD = $ws.Cells.Item($row, 4)
    $cD.Value2 = $r.D
    $cD.Font.Name = "Calibri"
    $cD.Font.Size = 10

    $cE = $ws.Cells.Item($row, 5)
    $cE.Value2 = $r.E
    $cE.Font.Name = "Calibri"
    $cE.Font.Size = 10

    $cF = $ws.Cells.Item($row, 6)
    $cF.Value2 = $r.F
    $cF.NumberFormat = "0.00%"
    $cF.Font.Name = "Calibri"
    $cF.Font.Size = 10

    $cG = $ws.Cells.Item($row, 7)
    $cG.Value2 = $r.G
    $cG.Font.Name = "Calibri"
    $cG.Font.Size = 10

    $cH = $ws.Cells.Item($row, 8)
    $cH.Value2 = $r.H
    $cH.NumberFormat = "dd/mm/yyyy"
    $cH.Font.Name = "Calibri"
    $cH.Font.Size = 10

    $cI = $ws.Cells.Item($row, 9)
    $cI.Value2 = $r.I
    $cI.NumberFormat = "$#,##0.00"
    $cI.Font.Name = "Calibri"
    $cI.Font.Size = 10
}

# Update selection to match the final saved view state
$ws.Range("A31:XFD31").Select()
